# Applies the "Updated cryptos list on Sun Aug 11 18:36:46 UTC 2024 with
# GitHub Actions" commit: refreshed Price (D) / Volume(1h) (E) figures for
# every coin row, plus a re-shuffle of the Bittensor / Stacks / Filecoin
# block (rows 40-42) with refreshed data for those three coins as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
  ,@(2, "D", "60.147.29")
  ,@(2, "E", "  -1.12%  ")
  ,@(3, "D", "2.623.18")
  ,@(3, "E", "  +0.91%  ")
  ,@(4, "E", "  +0.00%  ")
  ,@(5, "D", "518.74")
  ,@(5, "E", "  -0.27%  ")
  ,@(6, "D", "147.59")
  ,@(6, "E", "  -4.40%  ")
  ,@(7, "D", "0.999")
  ,@(7, "E", "  +0.06%  ")
  ,@(8, "D", "0.570")
  ,@(8, "E", "  -3.60%  ")
  ,@(9, "D", "2.628.06")
  ,@(9, "E", "  +0.67%  ")
  ,@(10, "D", "6.31")
  ,@(10, "E", "  -5.48%  ")
  ,@(11, "D", "0.105")
  ,@(11, "E", "  -0.49%  ")
  ,@(12, "D", "0.340")
  ,@(12, "E", "  -2.22%  ")
  ,@(13, "E", "  -0.68%  ")
  ,@(14, "D", "3.087.16")
  ,@(14, "E", "  +0.96%  ")
  ,@(15, "D", "60.149.53")
  ,@(15, "E", "  -1.17%  ")
  ,@(16, "D", "21.15")
  ,@(16, "E", "  -2.37%  ")
  ,@(17, "D", "0.0000138")
  ,@(17, "E", "  -1.82%  ")
  ,@(18, "D", "2.630.31")
  ,@(18, "E", "  +0.77%  ")
  ,@(19, "D", "4.62")
  ,@(19, "E", "  -2.26%  ")
  ,@(20, "D", "341.28")
  ,@(20, "E", "  -3.22%  ")
  ,@(21, "D", "10.40")
  ,@(21, "E", "  -1.42%  ")
  ,@(22, "D", "6.10")
  ,@(22, "E", "  -1.70%  ")
  ,@(23, "D", "0.995")
  ,@(23, "E", "  -0.44%  ")
  ,@(24, "D", "61.12")
  ,@(24, "E", "  +0.12%  ")
  ,@(25, "D", "0.417")
  ,@(25, "E", "  -2.14%  ")
  ,@(26, "D", "0.999")
  ,@(26, "E", "  +0.11%  ")
  ,@(27, "E", "  -3.54%  ")
  ,@(28, "D", "0.0₃0807")
  ,@(28, "E", "  -4.58%  ")
  ,@(29, "D", "7.02")
  ,@(29, "E", "  -4.45%  ")
  ,@(30, "D", "1.00")
  ,@(30, "E", "  +0.00%  ")
  ,@(31, "E", "  -0.89%  ")
  ,@(32, "D", "5.95")
  ,@(32, "E", "  -5.32%  ")
  ,@(33, "D", "18.89")
  ,@(33, "E", "  -2.39%  ")
  ,@(34, "D", "150.23")
  ,@(34, "E", "  +0.80%  ")
  ,@(35, "D", "3.93")
  ,@(35, "E", "  -6.62%  ")
  ,@(36, "E", "  -1.53%  ")
  ,@(37, "D", "1.13")
  ,@(37, "E", "  -5.54%  ")
  ,@(38, "D", "0.855")
  ,@(38, "E", "  +1.07%  ")
  ,@(39, "D", "36.67")
  ,@(39, "E", "  +0.64%  ")
  ,@(40, "B", "Stacks")
  ,@(40, "C", "https://coinranking.com/coin/mMPrMcB7+stacks-stx")
  ,@(40, "D", "1.42")
  ,@(40, "E", "  -4.71%  ")
  ,@(41, "B", "Filecoin")
  ,@(41, "C", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil")
  ,@(41, "D", "3.63")
  ,@(41, "E", "  -4.12%  ")
  ,@(42, "B", "Bittensor")
  ,@(42, "C", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao")
  ,@(42, "D", "290.64")
  ,@(42, "E", "  +1.48%  ")
  ,@(43, "D", "0.627")
  ,@(43, "E", "  +0.39%  ")
  ,@(44, "D", "0.0998")
  ,@(44, "E", "  -1.16%  ")
  ,@(45, "E", "  +0.11%  ")
  ,@(46, "D", "0.0546")
  ,@(46, "E", "  -2.48%  ")
  ,@(47, "D", "19.39")
  ,@(47, "E", "  -0.88%  ")
  ,@(48, "E", "  +0.79%  ")
  ,@(49, "E", "  -2.14%  ")
  ,@(50, "D", "4.64")
  ,@(50, "E", "  -4.37%  ")
  ,@(51, "D", "1.955.86")
  ,@(51, "E", "  -0.20%  ")
)

foreach ($item in $changes) {
  $row = $item[0]
  $col = $item[1]
  $val = $item[2]
  $addr = "$col$row"
  $rng = $ws.Range($addr)
  # Force text storage so numeric-looking values (e.g. "0.570", "1.00",
  # "60.147.29") keep their exact original formatting instead of being
  # auto-coerced into numbers by Excel's input parser.
  $rng.NumberFormat = "@"
  $rng.Value = $val
  # Drop back to the default (unstyled) cell style so no stray style index
  # is left on the cell - these cells carry no explicit style in the source.
  $rng.Style = "Normal"
}

Write-Output "Applied $($changes.Count) cell updates"
